# Refresh the crypto price/volume snapshot (GitHub Actions scheduled update).
# For numeric-looking "Price" strings (column D) we force NumberFormat="@"
# (Text) before assigning .Value, so Excel stores the literal text (e.g.
# "1.00", "0.0233") instead of silently coercing it to a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.499.66"
$ws.Range("E2").Value = "  -3.22%  "

$ws.Range("D3").Value = "2.906.97"
$ws.Range("E3").Value = "  -3.14%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "526.92"
$ws.Range("E5").Value = "  -4.99%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.96"
$ws.Range("E6").Value = "  -6.56%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.14%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.549"
$ws.Range("E8").Value = "  -3.51%  "

$ws.Range("D9").Value = "2.914.46"
$ws.Range("E9").Value = "  -2.74%  "

$ws.Range("E10").Value = "  -4.94%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.89"
$ws.Range("E11").Value = "  -6.24%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.356"
$ws.Range("E12").Value = "  -2.96%  "

$ws.Range("D13").Value = "3.411.56"
$ws.Range("E13").Value = "  -3.22%  "

$ws.Range("E14").Value = "  +1.12%  "

$ws.Range("D15").Value = "60.590.99"
$ws.Range("E15").Value = "  -3.20%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.64"
$ws.Range("E16").Value = "  -4.52%  "

$ws.Range("D17").Value = "2.906.34"
$ws.Range("E17").Value = "  -3.38%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0000140"
$ws.Range("E18").Value = "  -6.03%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.93"
$ws.Range("E19").Value = "  -3.47%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.51"
$ws.Range("E20").Value = "  -3.60%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "360.01"
$ws.Range("E21").Value = "  -8.63%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.55"
$ws.Range("E22").Value = "  -1.91%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.04%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.12"
$ws.Range("E24").Value = "  -3.11%  "

$ws.Range("D25").Value = "3.010.25"
$ws.Range("E25").Value = "  -4.30%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.447"
$ws.Range("E26").Value = "  -4.15%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.182"
$ws.Range("E27").Value = "  -2.36%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  -0.23%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.79"
$ws.Range("E29").Value = "  -8.86%  "

$ws.Range("D30").Value = "0.0₃0852"
$ws.Range("E30").Value = "  -12.24%  "

$ws.Range("E31").Value = "  -0.08%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.66"
$ws.Range("E32").Value = "  -4.99%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.43"
$ws.Range("E33").Value = "  -5.48%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "150.04"
$ws.Range("E34").Value = "  -6.43%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.31"
$ws.Range("E35").Value = "  -7.78%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.53"
$ws.Range("E36").Value = "  -8.25%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.985"
$ws.Range("E37").Value = "  -9.38%  "

$ws.Range("E38").Value = "  -7.55%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "38.00"
$ws.Range("E39").Value = "  +1.21%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.48"
$ws.Range("E40").Value = "  -6.29%  "

$ws.Range("D41").Value = "2.331.58"
$ws.Range("E41").Value = "  -5.40%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.645"
$ws.Range("E42").Value = "  -2.27%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.64"
$ws.Range("E43").Value = "  -6.96%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "20.68"
$ws.Range("E44").Value = "  -8.08%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0569"
$ws.Range("E45").Value = "  -4.52%  "

$ws.Range("E46").Value = "  +0.00%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.91"
$ws.Range("E47").Value = "  +0.42%  "

$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0233"
$ws.Range("E48").Value = "  -5.96%  "

$ws.Range("B49").Value = "WhiteBITCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.35"
$ws.Range("E49").Value = "  -1.46%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0926"
$ws.Range("E50").Value = "  -2.81%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "250.38"
$ws.Range("E51").Value = "  -4.54%  "
